# chore: update Sheets via scheduled runner
# Refreshes cached marketboard-derived profit figures (currentAveragePrice*,
# LevePrice*, LeveProfit*) across all job sheets.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 257.375
$ws.Range("I2").Value = 220.66667
$ws.Range("J2").Value = 367.5
$ws.Range("K2").Value = 220.66667
$ws.Range("L2").Value = 367.5
$ws.Range("M2").Value = -107.66667
$ws.Range("N2").Value = -593.5
$ws.Range("H11").Value = 195.26666
$ws.Range("I11").Value = 195.26666
$ws.Range("K11").Value = 195.26666
$ws.Range("M11").Value = -55.26666
$ws.Range("H92").Value = 3337.05
$ws.Range("I92").Value = 1628.4
$ws.Range("K92").Value = 1628.4
$ws.Range("M92").Value = -380.4000000000001
$ws.Range("H98").Value = 231606.89
$ws.Range("I98").Value = 898.5217
$ws.Range("K98").Value = 898.5217
$ws.Range("M98").Value = 599.4783
$ws.Range("H103").Value = 1079.5
$ws.Range("J103").Value = 1399.4
$ws.Range("L103").Value = 4198.200000000001
$ws.Range("N103").Value = -5370.200000000001
$ws.Range("H113").Value = 5239.6
$ws.Range("J113").Value = 5239.6
$ws.Range("L113").Value = 5239.6
$ws.Range("N113").Value = -11747.6
$ws.Range("H122").Value = 231606.89
$ws.Range("I122").Value = 898.5217
$ws.Range("K122").Value = 2695.5651
$ws.Range("M122").Value = -245.5650999999998
$ws.Range("H138").Value = 2852.447
$ws.Range("I138").Value = 2026.6666
$ws.Range("J138").Value = 3029.4
$ws.Range("K138").Value = 6079.9998
$ws.Range("L138").Value = 9088.200000000001
$ws.Range("M138").Value = -939.9997999999996
$ws.Range("N138").Value = -19368.2

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 6061.3887
$ws.Range("I2").Value = 397
$ws.Range("J2").Value = 34383.332
$ws.Range("K2").Value = 397
$ws.Range("L2").Value = 34383.332
$ws.Range("M2").Value = -284
$ws.Range("N2").Value = -34609.332
$ws.Range("H32").Value = 7612.7144
$ws.Range("I32").Value = 5992.1763
$ws.Range("K32").Value = 5992.1763
$ws.Range("M32").Value = -5705.1763
$ws.Range("H97").Value = 533.7368
$ws.Range("I97").Value = 424.5
$ws.Range("J97").Value = 2500
$ws.Range("K97").Value = 424.5
$ws.Range("L97").Value = 2500
$ws.Range("M97").Value = 71.5
$ws.Range("N97").Value = -3492
$ws.Range("H102").Value = 3838.75
$ws.Range("I102").Value = 3475
$ws.Range("J102").Value = 3960
$ws.Range("K102").Value = 3475
$ws.Range("L102").Value = 3960
$ws.Range("M102").Value = -1853
$ws.Range("N102").Value = -7204
$ws.Range("H116").Value = 6061.3887
$ws.Range("I116").Value = 397
$ws.Range("J116").Value = 34383.332
$ws.Range("K116").Value = 397
$ws.Range("L116").Value = 34383.332
$ws.Range("M116").Value = 1897
$ws.Range("N116").Value = -38971.332
$ws.Range("H122").Value = 2999.9656
$ws.Range("I122").Value = 2379.353
$ws.Range("K122").Value = 7138.059
$ws.Range("M122").Value = -4688.059
$ws.Range("H132").Value = 2460.1035
$ws.Range("I132").Value = 1752.9259
$ws.Range("K132").Value = 5258.7777
$ws.Range("M132").Value = -2728.7777
$ws.Range("H133").Value = 75087
$ws.Range("J133").Value = 75087
$ws.Range("L133").Value = 75087
$ws.Range("N133").Value = -80147

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 6061.3887
$ws.Range("I3").Value = 397
$ws.Range("J3").Value = 34383.332
$ws.Range("K3").Value = 397
$ws.Range("L3").Value = 34383.332
$ws.Range("M3").Value = -283
$ws.Range("N3").Value = -34611.332
$ws.Range("H99").Value = 1839.5
$ws.Range("I99").Value = 1562.1111
$ws.Range("K99").Value = 1562.1111
$ws.Range("M99").Value = -64.11110000000008
$ws.Range("H123").Value = 48302.855
$ws.Range("J123").Value = 48302.855
$ws.Range("L123").Value = 48302.855
$ws.Range("N123").Value = -58102.855
$ws.Range("H134").Value = 1848.6086
$ws.Range("I134").Value = 1965.25
$ws.Range("K134").Value = 5895.75
$ws.Range("M134").Value = -3360.75

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 147.10527
$ws.Range("I7").Value = 41.22222
$ws.Range("J7").Value = 242.4
$ws.Range("K7").Value = 41.22222
$ws.Range("L7").Value = 242.4
$ws.Range("M7").Value = 71.77778000000001
$ws.Range("N7").Value = -468.4
$ws.Range("H32").Value = 3477
$ws.Range("I32").Value = 3954
$ws.Range("J32").Value = 3000
$ws.Range("K32").Value = 3954
$ws.Range("L32").Value = 3000
$ws.Range("M32").Value = -3638
$ws.Range("N32").Value = -3632
$ws.Range("H95").Value = 20109.875
$ws.Range("J95").Value = 20109.875
$ws.Range("L95").Value = 20109.875
$ws.Range("N95").Value = -25601.875
$ws.Range("H96").Value = 24168
$ws.Range("J96").Value = 24168
$ws.Range("L96").Value = 24168
$ws.Range("N96").Value = -29660
$ws.Range("H99").Value = 3301.1667
$ws.Range("I99").Value = 2834.8572
$ws.Range("K99").Value = 2834.8572
$ws.Range("M99").Value = -1336.8572
$ws.Range("H122").Value = 4737
$ws.Range("I122").Value = 2032.1333
$ws.Range("J122").Value = 11499.167
$ws.Range("K122").Value = 6096.3999
$ws.Range("L122").Value = 34497.501
$ws.Range("M122").Value = -3646.3999
$ws.Range("N122").Value = -39397.501
$ws.Range("H126").Value = 3301.1667
$ws.Range("I126").Value = 2834.8572
$ws.Range("K126").Value = 8504.571599999999
$ws.Range("M126").Value = -6034.571599999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H136").Value = 47622604
$ws.Range("I136").Value = 66669644
$ws.Range("K136").Value = 200008932
$ws.Range("M136").Value = -200003832
$ws.Range("H141").Value = 7026.533
$ws.Range("I141").Value = 3786.5
$ws.Range("J141").Value = 13506.6
$ws.Range("K141").Value = 11359.5
$ws.Range("L141").Value = 40519.8
$ws.Range("M141").Value = -6179.5
$ws.Range("N141").Value = -50879.8

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H5").Value = 49.5
$ws.Range("I5").Value = 49.5
$ws.Range("K5").Value = 49.5
$ws.Range("M5").Value = 62.5
$ws.Range("H107").Value = 405.72
$ws.Range("I107").Value = 405.375
$ws.Range("J107").Value = 406.33334
$ws.Range("K107").Value = 405.375
$ws.Range("L107").Value = 406.33334
$ws.Range("M107").Value = 1514.625
$ws.Range("N107").Value = -4246.33334
$ws.Range("H113").Value = 2724.84
$ws.Range("J113").Value = 3394.2
$ws.Range("L113").Value = 3394.2
$ws.Range("N113").Value = -7734.2
$ws.Range("H126").Value = 4799.5884
$ws.Range("I126").Value = 3198.6667
$ws.Range("K126").Value = 9596.000100000001
$ws.Range("M126").Value = -7126.000100000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 17655.555
$ws.Range("I132").Value = 16980.834
$ws.Range("K132").Value = 50942.50199999999
$ws.Range("M132").Value = -48412.50199999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H82").Value = 49998.668
$ws.Range("J82").Value = 49998
$ws.Range("L82").Value = 49998
$ws.Range("N82").Value = -50764
$ws.Range("H85").Value = 49998.668
$ws.Range("J85").Value = 49998
$ws.Range("L85").Value = 49998
$ws.Range("N85").Value = -52650
